$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.801.97"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "3.467.05"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'414.54"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").Value = "'129.85"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("D7").Value = "'0.627"
$ws.Range("E7").Value = "  -0.75%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.725"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("D10").Value = "'0.153"
$ws.Range("E10").Value = "  +8.96%  "
$ws.Range("D11").Value = "'42.47"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").Value = "'0.0000229"
$ws.Range("E12").Value = "  +6.23%  "
$ws.Range("E13").Value = "  +5.65%  "
$ws.Range("D14").Value = "4.019.65"
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("D15").Value = "'0.140"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").Value = "'20.49"
$ws.Range("E16").Value = "  -3.35%  "
$ws.Range("D17").Value = "3.462.03"
$ws.Range("E17").Value = "  +1.32%  "
$ws.Range("D18").Value = "'12.59"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("D19").Value = "'1.07"
$ws.Range("E19").Value = "  -1.27%  "
$ws.Range("D20").Value = "62.744.91"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").Value = "'459.34"
$ws.Range("E21").Value = "  +1.98%  "
$ws.Range("D22").Value = "'90.33"
$ws.Range("E22").Value = "  -1.43%  "
$ws.Range("D23").Value = "'3.27"
$ws.Range("E23").Value = "  +1.90%  "
$ws.Range("D24").Value = "'13.21"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("D25").Value = "'10.71"
$ws.Range("E25").Value = "  +14.62%  "
$ws.Range("D26").Value = "'3.31"
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("D27").Value = "'33.28"
$ws.Range("E27").Value = "  +0.88%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("D30").Value = "'12.05"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("E32").Value = "  -1.38%  "
$ws.Range("E33").Value = "  -2.01%  "
$ws.Range("D34").Value = "'40.78"
$ws.Range("E34").Value = "  -4.37%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'58.29"
$ws.Range("E36").Value = "  +8.24%  "
$ws.Range("D37").Value = "'0.0490"
$ws.Range("E37").Value = "  -2.37%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'3.07"
$ws.Range("E39").Value = "  +3.91%  "
$ws.Range("D40").Value = "'149.90"
$ws.Range("E40").Value = "  +5.07%  "
$ws.Range("D41").Value = "'2.71"
$ws.Range("E41").Value = "  +6.08%  "
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("D43").Value = "'0.320"
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("D45").Value = "'4.41"
$ws.Range("E45").Value = "  +3.68%  "
$ws.Range("D46").Value = "'2.06"
$ws.Range("E46").Value = "  +2.85%  "
$ws.Range("D47").Value = "'2.38"
$ws.Range("E47").Value = "  +11.61%  "
$ws.Range("D48").Value = "0.0₃0561"
$ws.Range("E48").Value = "  +33.97%  "
$ws.Range("D49").Value = "'16.33"
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("D50").Value = "'22.22"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("D51").Value = "'0.140"
$ws.Range("E51").Value = "  -4.12%  "
